# Data Dictionary.xlsx update
#
# The "Villains_Abilities" and "Heroes_Abilities" tables each had a "Level"
# field row that is removed, and their trailing "Damage" field is renamed to
# "Overall Damage".
#
# Before (Villains_Abilities, rows 13-16 / Heroes_Abilities, rows 17-20):
#   Villain_ID | int   | PK / FK
#   Ability_ID | int   | PK / FK
#   Level      | int   |
#   Damage     | float |
#
# After (Villains_Abilities, rows 13-15 / Heroes_Abilities, rows 16-18):
#   Villain_ID | int   | PK / FK
#   Ability_ID | int   | PK / FK
#   Overall Damage | float |

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Level" row of the second (lower) table first so the row
# numbers of the first table are unaffected while we still need them.
$ws.Rows(19).Delete()
$ws.Rows(15).Delete()

# Rename the remaining "Damage" rows (now the last row of each shrunk
# table) to "Overall Damage".
$ws.Range("B15").Value = "Overall Damage"
$ws.Range("B18").Value = "Overall Damage"

# Match the author's final selection/active cell.
$ws.Range("D19").Select()
